# Update the "Submission deadline" date on slide 1 (Rounded Rectangle 26)
# from "Midnight 1st Nov 2023" to "Midnight 20th Oct 2023", keeping the
# bold formatting throughout and the superscript ("th") formatting that
# previously applied to "st".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Rounded Rectangle 26")
$tf = $sh.TextFrame
$tr = $tf.TextRange

$fullText = $tr.Text
$oldPhrase = "Midnight 1st Nov 2023"
$newPhrase = "Midnight 20th Oct 2023"

$idx = $fullText.IndexOf($oldPhrase)
if ($idx -lt 0) {
    throw "Could not find '" + $oldPhrase + "' in shape text."
}
$base = $idx + 1

# Replace the whole phrase first (single run, inherits the bold formatting
# already present on "Midnight 1st Nov 2023").
$whole = $tr.Characters($base, $oldPhrase.Length)
$whole.Text = $newPhrase

# Now re-split the replaced phrase into the individual runs so each piece
# carries the right formatting:
#   "Midnight " -> bold
#   "20"        -> bold
#   "th"        -> bold + superscript
#   " Oct "     -> bold
#   "2023"      -> bold

$r1 = $tr.Characters($base + 0, 9)   # "Midnight "
$r1.Font.Bold = -1

$r2 = $tr.Characters($base + 9, 2)   # "20"
$r2.Font.Bold = -1

$r3 = $tr.Characters($base + 11, 2)  # "th"
$r3.Font.Bold = -1
$r3.Font.Superscript = -1

$r4 = $tr.Characters($base + 13, 5)  # " Oct "
$r4.Font.Bold = -1

$r5 = $tr.Characters($base + 18, 4)  # "2023"
$r5.Font.Bold = -1
